# Add 2022-Q3 data:
#   1. "总计" (summary) sheet: insert a new data row for 2022-Q3 at the top of
#      the data block (row 2), pushing the existing quarters down by one row.
#   2. Insert a brand-new worksheet named "2022-Q3" right after "总计" and
#      before "2022-Q2", populated with the per-fund holding breakdown for
#      that quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "总计" sheet — shift existing rows down and insert the new 2022-Q3 row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Snapshot the existing 7 data rows (rows 2-8 => B,C,D) before overwriting.
$existing = @()
for ($r = 2; $r -le 8; $r++) {
    $existing += ,@($summary.Cells.Item($r,2).Value2, $summary.Cells.Item($r,3).Value2, $summary.Cells.Item($r,4).Value2)
}

# Re-write them one row lower (row 2..8 -> row 3..9), bumping the index in
# column A by one to match.
for ($i = 0; $i -lt $existing.Length; $i++) {
    $newRow = $i + 3
    $summary.Cells.Item($newRow,1).Value = $i + 1
    $summary.Cells.Item($newRow,2).Value = $existing[$i][0]
    $summary.Cells.Item($newRow,3).Value = $existing[$i][1]
    $summary.Cells.Item($newRow,4).Value = $existing[$i][2]
}

# Row 9 is brand new territory (sheet used to stop at row 8) - copy the
# number format / border / bold styling of column A down onto it.
$summary.Range("A8").Copy()
$summary.Range("A9").PasteSpecial(-4122)

# Write the new 2022-Q3 row at the top of the data block.
$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q3"
$summary.Cells.Item(2,3).Value = 14
$summary.Cells.Item(2,4).Value = 2.67

# ---------------------------------------------------------------------------
# 2. Insert the new "2022-Q3" worksheet before "2022-Q2" (i.e. right after
#    "总计"). Duplicate "2022-Q2" itself so the new sheet inherits identical
#    sheet-level formatting (column headers, styles, page setup, ...), then
#    trim/overwrite the data rows.
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item(2)
$q2.Copy($q2, $null)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# 2022-Q2 has 38 data rows (rows 2-39); 2022-Q3 only needs 14 (rows 2-15).
# Drop the now-unused trailing rows entirely.
$q3.Range("A16:A39").EntireRow.Delete()

$rows = @(
    @("005739", "富国转型机遇混合", "44.36", "81.74", "2.86", "1.2687", 9),
    @("014057", "富国金安均衡精选混合A", "21.31", "77.70", "2.81", "0.5988", 8),
    @("006527", "富国优质发展混合A", "15.53", "78.90", "2.66", "0.4131", 10),
    @("006528", "富国优质发展混合C", "3.64", "78.90", "2.66", "0.0968", 10),
    @("004895", "华商鑫安灵活配置混合", "2.11", "92.54", "4.27", "0.0901", 5),
    @("014058", "富国金安均衡精选混合C", "1.82", "77.70", "2.81", "0.0511", 8),
    @("310368", "申万菱信竞争优势混合A", "1.05", "92.76", "4.74", "0.0498", 5),
    @("519677", "银河定投宝腾讯济安指数", "2.88", "92.40", "1.54", "0.0444", 3),
    @("005009", "申万菱信行业轮动股票A", "0.62", "92.87", "4.76", "0.0295", 2),
    @("008116", "银华沪深股通精选混合", "0.50", "88.46", "3.94", "0.0197", 9),
    @("015173", "申万菱信竞争优势混合C", "0.11", "92.76", "4.74", "0.0052", 5),
    @("014692", "中加量化研选混合型证券投资基金C", "0.19", "68.41", "2.37", "0.0045", 2),
    @("015157", "申万菱信行业轮动股票C", "0.04", "92.87", "4.76", "0.0019", 2),
    @("014691", "中加量化研选混合型证券投资基金A", "0.05", "68.41", "2.37", "0.0012", 2)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    $q3.Cells.Item($r, 1).Value = $i

    # Columns B..G must stay text (fund codes have leading zeros, the
    # numeric-looking figures keep trailing zeros) - force text with a
    # leading apostrophe, then strip the resulting quote-prefix style back
    # to the sheet's untouched default.
    for ($c = 0; $c -lt 6; $c++) {
        $cell = $q3.Cells.Item($r, $c + 2)
        $cell.Value = "'" + $row[$c]
        $cell.ClearFormats()
    }

    # Column H (rank) is a genuine number.
    $q3.Cells.Item($r, 8).Value = $row[6]
}

# Restore the original active sheet/tab selection.
$summary.Activate()

Write-Host "2022-Q3 sheet inserted"
